$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2249.5833  # H43: 1875 -> 2249.5833
$ws.Cells.Item(43, 9).Value = 2200  # I43: 2000 -> 2200
$ws.Cells.Item(43, 10).Value = 2285  # J43: 1750 -> 2285
$ws.Cells.Item(43, 11).Value = 2200  # K43: 2000 -> 2200
$ws.Cells.Item(43, 12).Value = 2285  # L43: 1750 -> 2285
$ws.Cells.Item(43, 13).Value = -2131  # M43: -1931 -> -2131
$ws.Cells.Item(43, 14).Value = -2423  # N43: -1888 -> -2423

$ws.Cells.Item(137, 8).Value = 2681.05  # H137: 2822.3157 -> 2681.05
$ws.Cells.Item(137, 9).Value = 1581  # I137: 1779 -> 1581
$ws.Cells.Item(137, 11).Value = 4743  # K137: 5337 -> 4743
$ws.Cells.Item(137, 13).Value = -2193  # M137: -2787 -> -2193

$ws.Cells.Item(138, 8).Value = 11641.929  # H138: 10268.23 -> 11641.929
$ws.Cells.Item(138, 9).Value = 999.75  # I138: 1099.2 -> 999.75
$ws.Cells.Item(138, 10).Value = 15898.8  # J138: 15998.875 -> 15898.8
$ws.Cells.Item(138, 11).Value = 2999.25  # K138: 3297.6 -> 2999.25
$ws.Cells.Item(138, 12).Value = 47696.39999999999  # L138: 47996.625 -> 47696.39999999999
$ws.Cells.Item(138, 13).Value = 2140.75  # M138: 1842.4 -> 2140.75
$ws.Cells.Item(138, 14).Value = -57976.39999999999  # N138: -58276.625 -> -57976.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18917.916  # H32: 19460 -> 18917.916
$ws.Cells.Item(32, 9).Value = 18546.818  # I32: 19138.182 -> 18546.818
$ws.Cells.Item(32, 11).Value = 18546.818  # K32: 19138.182 -> 18546.818
$ws.Cells.Item(32, 13).Value = -18259.818  # M32: -18851.182 -> -18259.818

$ws.Cells.Item(63, 8).Value = 7033.3335  # H63: 5552.75 -> 7033.3335
$ws.Cells.Item(63, 9).Value = 2903.3333  # I63: 2988.6667 -> 2903.3333
$ws.Cells.Item(63, 10).Value = 11163.333  # J63: 13245 -> 11163.333
$ws.Cells.Item(63, 11).Value = 2903.3333  # K63: 2988.6667 -> 2903.3333
$ws.Cells.Item(63, 12).Value = 11163.333  # L63: 13245 -> 11163.333
$ws.Cells.Item(63, 13).Value = -2217.3333  # M63: -2302.6667 -> -2217.3333
$ws.Cells.Item(63, 14).Value = -12535.333  # N63: -14617 -> -12535.333

$ws.Cells.Item(66, 8).Value = 7033.3335  # H66: 5552.75 -> 7033.3335
$ws.Cells.Item(66, 9).Value = 2903.3333  # I66: 2988.6667 -> 2903.3333
$ws.Cells.Item(66, 10).Value = 11163.333  # J66: 13245 -> 11163.333
$ws.Cells.Item(66, 11).Value = 14516.6665  # K66: 14943.3335 -> 14516.6665
$ws.Cells.Item(66, 12).Value = 55816.665  # L66: 66225 -> 55816.665
$ws.Cells.Item(66, 13).Value = -11084.6665  # M66: -11511.3335 -> -11084.6665
$ws.Cells.Item(66, 14).Value = -62680.665  # N66: -73089 -> -62680.665

$ws.Cells.Item(132, 8).Value = 3399.5715  # H132: 3524.6667 -> 3399.5715
$ws.Cells.Item(132, 9).Value = 2933  # I132: 3049.6667 -> 2933
$ws.Cells.Item(132, 10).Value = 3749.5  # J132: 3999.6667 -> 3749.5
$ws.Cells.Item(132, 11).Value = 8799  # K132: 9149.000100000001 -> 8799
$ws.Cells.Item(132, 12).Value = 11248.5  # L132: 11999.0001 -> 11248.5
$ws.Cells.Item(132, 13).Value = -6269  # M132: -6619.000100000001 -> -6269
$ws.Cells.Item(132, 14).Value = -16308.5  # N132: -17059.0001 -> -16308.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2500  # H134: 4000 -> 2500
$ws.Cells.Item(134, 9).Value = 2500  # I134: 4000 -> 2500
$ws.Cells.Item(134, 11).Value = 7500  # K134: 12000 -> 7500
$ws.Cells.Item(134, 13).Value = -4965  # M134: -9465 -> -4965

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2010.6471  # H31: 2146.0667 -> 2010.6471
$ws.Cells.Item(31, 9).Value = 1691.6428  # I31: 1807.75 -> 1691.6428
$ws.Cells.Item(31, 11).Value = 1691.6428  # K31: 1807.75 -> 1691.6428
$ws.Cells.Item(31, 13).Value = -1396.6428  # M31: -1512.75 -> -1396.6428

$ws.Cells.Item(34, 8).Value = 2010.6471  # H34: 2146.0667 -> 2010.6471
$ws.Cells.Item(34, 9).Value = 1691.6428  # I34: 1807.75 -> 1691.6428
$ws.Cells.Item(34, 11).Value = 1691.6428  # K34: 1807.75 -> 1691.6428
$ws.Cells.Item(34, 13).Value = -1489.6428  # M34: -1605.75 -> -1489.6428

$ws.Cells.Item(58, 8).Value = 3454.6  # H58: 3538.7778 -> 3454.6
$ws.Cells.Item(58, 9).Value = 3698.4285  # I58: 3865 -> 3698.4285
$ws.Cells.Item(58, 10).Value = 2885.6667  # J58: 2886.3333 -> 2885.6667
$ws.Cells.Item(58, 11).Value = 3698.4285  # K58: 3865 -> 3698.4285
$ws.Cells.Item(58, 12).Value = 2885.6667  # L58: 2886.3333 -> 2885.6667
$ws.Cells.Item(58, 13).Value = -3495.4285  # M58: -3662 -> -3495.4285
$ws.Cells.Item(58, 14).Value = -3291.6667  # N58: -3292.3333 -> -3291.6667

$ws.Cells.Item(86, 8).Value = 19496.5  # H86: 27748.75 -> 19496.5
$ws.Cells.Item(86, 9).Value = 27746  # I86: 52500 -> 27746
$ws.Cells.Item(86, 11).Value = 27746  # K86: 52500 -> 27746
$ws.Cells.Item(86, 13).Value = -26623  # M86: -51377 -> -26623

$ws.Cells.Item(89, 8).Value = 19496.5  # H89: 27748.75 -> 19496.5
$ws.Cells.Item(89, 9).Value = 27746  # I89: 52500 -> 27746
$ws.Cells.Item(89, 11).Value = 138730  # K89: 262500 -> 138730
$ws.Cells.Item(89, 13).Value = -133114  # M89: -256884 -> -133114

$ws.Cells.Item(94, 8).Value = 4999.3335  # H94: 4999.6665 -> 4999.3335
$ws.Cells.Item(94, 9).Value = 4999.3335  # I94: 5000 -> 4999.3335
$ws.Cells.Item(94, 10).Value = 0  # J94: 4999 -> 0
$ws.Cells.Item(94, 11).Value = 4999.3335  # K94: 5000 -> 4999.3335
$ws.Cells.Item(94, 12).Value = 0  # L94: 4999 -> 0
$ws.Cells.Item(94, 13).Value = -4548.3335  # M94: -4549 -> -4548.3335
$ws.Cells.Item(94, 14).ClearContents()  # N94: -5901 -> (removed)

$ws.Cells.Item(122, 8).Value = 1218.6  # H122: 1115.4166 -> 1218.6
$ws.Cells.Item(122, 9).Value = 1347  # I122: 1219.5555 -> 1347
$ws.Cells.Item(122, 10).Value = 705  # J122: 803 -> 705
$ws.Cells.Item(122, 11).Value = 4041  # K122: 3658.6665 -> 4041
$ws.Cells.Item(122, 12).Value = 2115  # L122: 2409 -> 2115
$ws.Cells.Item(122, 13).Value = -1591  # M122: -1208.6665 -> -1591
$ws.Cells.Item(122, 14).Value = -7015  # N122: -7309 -> -7015

$ws.Cells.Item(136, 8).Value = 3454.6  # H136: 3538.7778 -> 3454.6
$ws.Cells.Item(136, 9).Value = 3698.4285  # I136: 3865 -> 3698.4285
$ws.Cells.Item(136, 10).Value = 2885.6667  # J136: 2886.3333 -> 2885.6667
$ws.Cells.Item(136, 11).Value = 11095.2855  # K136: 11595 -> 11095.2855
$ws.Cells.Item(136, 12).Value = 8657.000100000001  # L136: 8658.999899999999 -> 8657.000100000001
$ws.Cells.Item(136, 13).Value = -8545.2855  # M136: -9045 -> -8545.2855
$ws.Cells.Item(136, 14).Value = -13757.0001  # N136: -13758.9999 -> -13757.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1108.5652  # H5: 1195.409 -> 1108.5652
$ws.Cells.Item(5, 9).Value = 899.8  # I5: 1079.9 -> 899.8
$ws.Cells.Item(5, 10).Value = 1269.1538  # J5: 1291.6666 -> 1269.1538
$ws.Cells.Item(5, 11).Value = 2699.4  # K5: 3239.7 -> 2699.4
$ws.Cells.Item(5, 12).Value = 3807.4614  # L5: 3874.9998 -> 3807.4614
$ws.Cells.Item(5, 13).Value = -2587.4  # M5: -3127.7 -> -2587.4
$ws.Cells.Item(5, 14).Value = -4031.4614  # N5: -4098.9998 -> -4031.4614

$ws.Cells.Item(12, 8).Value = 95  # H12: 96.90909000000001 -> 95
$ws.Cells.Item(12, 10).Value = 95.71429000000001  # J12: 98.71429000000001 -> 95.71429000000001
$ws.Cells.Item(12, 12).Value = 287.14287  # L12: 296.14287 -> 287.14287
$ws.Cells.Item(12, 14).Value = -633.14287  # N12: -642.14287 -> -633.14287

$ws.Cells.Item(135, 8).Value = 1108.5652  # H135: 1195.409 -> 1108.5652
$ws.Cells.Item(135, 9).Value = 899.8  # I135: 1079.9 -> 899.8
$ws.Cells.Item(135, 10).Value = 1269.1538  # J135: 1291.6666 -> 1269.1538
$ws.Cells.Item(135, 11).Value = 8098.2  # K135: 9719.1 -> 8098.2
$ws.Cells.Item(135, 12).Value = 11422.3842  # L135: 11624.9994 -> 11422.3842
$ws.Cells.Item(135, 13).Value = -5563.2  # M135: -7184.1 -> -5563.2
$ws.Cells.Item(135, 14).Value = -16492.3842  # N135: -16694.9994 -> -16492.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 777  # H102: 0 -> 777
$ws.Cells.Item(102, 9).Value = 777  # I102: 0 -> 777
$ws.Cells.Item(102, 11).Value = 777  # K102: 0 -> 777
$ws.Cells.Item(102, 13).Value = 845  # M102: (new) -> 845

$ws.Cells.Item(126, 8).Value = 3998.5  # H126: 3993 -> 3998.5
$ws.Cells.Item(126, 9).Value = 4000  # I126: 3989.6667 -> 4000
$ws.Cells.Item(126, 10).Value = 3997  # J126: 3998 -> 3997
$ws.Cells.Item(126, 11).Value = 12000  # K126: 11969.0001 -> 12000
$ws.Cells.Item(126, 12).Value = 11991  # L126: 11994 -> 11991
$ws.Cells.Item(126, 13).Value = -9530  # M126: -9499.000100000001 -> -9530
$ws.Cells.Item(126, 14).Value = -16931  # N126: -16934 -> -16931

$ws.Cells.Item(132, 8).Value = 5249.6665  # H132: 3099.6365 -> 5249.6665
$ws.Cells.Item(132, 9).Value = 4000  # I132: 1514 -> 4000
$ws.Cells.Item(132, 11).Value = 12000  # K132: 4542 -> 12000
$ws.Cells.Item(132, 13).Value = -9470  # M132: -2012 -> -9470

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6964.25  # H22: 7108.7617 -> 6964.25
$ws.Cells.Item(22, 10).Value = 7385.7144  # J22: 7712.375 -> 7385.7144
$ws.Cells.Item(22, 12).Value = 7385.7144  # L22: 7712.375 -> 7385.7144
$ws.Cells.Item(22, 14).Value = -7975.7144  # N22: -8302.375 -> -7975.7144

$ws.Cells.Item(27, 8).Value = 6964.25  # H27: 7108.7617 -> 6964.25
$ws.Cells.Item(27, 10).Value = 7385.7144  # J27: 7712.375 -> 7385.7144
$ws.Cells.Item(27, 12).Value = 7385.7144  # L27: 7712.375 -> 7385.7144
$ws.Cells.Item(27, 14).Value = -7599.7144  # N27: -7926.375 -> -7599.7144

$ws.Cells.Item(40, 8).Value = 1000  # H40: 4400 -> 1000
$ws.Cells.Item(40, 9).Value = 1000  # I40: 4400 -> 1000
$ws.Cells.Item(40, 11).Value = 1000  # K40: 4400 -> 1000
$ws.Cells.Item(40, 13).Value = -864  # M40: -4264 -> -864

$ws.Cells.Item(82, 8).Value = 3876.6924  # H82: 3876.7693 -> 3876.6924
$ws.Cells.Item(82, 10).Value = 9123.25  # J82: 9123.5 -> 9123.25
$ws.Cells.Item(82, 12).Value = 9123.25  # L82: 9123.5 -> 9123.25
$ws.Cells.Item(82, 14).Value = -9845.25  # N82: -9845.5 -> -9845.25

$ws.Cells.Item(85, 8).Value = 3876.6924  # H85: 3876.7693 -> 3876.6924
$ws.Cells.Item(85, 10).Value = 9123.25  # J85: 9123.5 -> 9123.25
$ws.Cells.Item(85, 12).Value = 9123.25  # L85: 9123.5 -> 9123.25
$ws.Cells.Item(85, 14).Value = -11619.25  # N85: -11619.5 -> -11619.25

$ws.Cells.Item(136, 8).Value = 7127.6665  # H136: 7456.125 -> 7127.6665
$ws.Cells.Item(136, 10).Value = 14500  # J136: 19500 -> 14500
$ws.Cells.Item(136, 12).Value = 43500  # L136: 58500 -> 43500
$ws.Cells.Item(136, 14).Value = -48600  # N136: -63600 -> -48600

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 25213.334  # H54: 26880.5 -> 25213.334
$ws.Cells.Item(54, 10).Value = 25213.334  # J54: 26880.5 -> 25213.334
$ws.Cells.Item(54, 12).Value = 25213.334  # L54: 26880.5 -> 25213.334
$ws.Cells.Item(54, 14).Value = -26253.334  # N54: -27920.5 -> -26253.334

$ws.Cells.Item(74, 8).Value = 36243.6  # H74: 37650.75 -> 36243.6
$ws.Cells.Item(74, 10).Value = 34054.75  # J74: 35201.332 -> 34054.75
$ws.Cells.Item(74, 12).Value = 34054.75  # L74: 35201.332 -> 34054.75
$ws.Cells.Item(74, 14).Value = -35926.75  # N74: -37073.332 -> -35926.75

$ws.Cells.Item(77, 8).Value = 36243.6  # H77: 37650.75 -> 36243.6
$ws.Cells.Item(77, 10).Value = 34054.75  # J77: 35201.332 -> 34054.75
$ws.Cells.Item(77, 12).Value = 102164.25  # L77: 105603.996 -> 102164.25
$ws.Cells.Item(77, 14).Value = -111524.25  # N77: -114963.996 -> -111524.25

$ws.Cells.Item(107, 8).Value = 662.6  # H107: 712.1818 -> 662.6
$ws.Cells.Item(107, 9).Value = 605.1429000000001  # I107: 679.5 -> 605.1429000000001
$ws.Cells.Item(107, 10).Value = 796.6667  # J107: 799.3333 -> 796.6667
$ws.Cells.Item(107, 11).Value = 1815.4287  # K107: 2038.5 -> 1815.4287
$ws.Cells.Item(107, 12).Value = 2390.0001  # L107: 2397.9999 -> 2390.0001
$ws.Cells.Item(107, 13).Value = 104.5712999999998  # M107: -118.5 -> 104.5712999999998
$ws.Cells.Item(107, 14).Value = -6230.0001  # N107: -6237.9999 -> -6230.0001

$ws.Cells.Item(136, 8).Value = 2922.8  # H136: 3553.5 -> 2922.8
$ws.Cells.Item(136, 9).Value = 2922.8  # I136: 3553.5 -> 2922.8
$ws.Cells.Item(136, 11).Value = 8768.400000000001  # K136: 10660.5 -> 8768.400000000001
$ws.Cells.Item(136, 13).Value = -6218.400000000001  # M136: -8110.5 -> -6218.400000000001
